$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("floor7_sw_lab",   "JUNOS Base OS boot [12.3R12.4]",                              "ex3300-24p", "10.9.106.30/23", "1", "38", "34"),
    @("PTSW1_Floor7",    "JUNOS Base OS boot [12.3R12.4]",                              "ex3300-48p", "10.9.106.37/23", "4", "66", "39"),
    @("PTSW2_Floor7",    "JUNOS Base OS boot [12.3R12.4]",                              "ex3300-48p", "10.9.106.38/23", "4", "71", "38"),
    @("oren-flr1sw-B1",  "JUNOS OS Kernel 32-bit  [20180119.e26d166_builder_master]",   "ex3400-48p", "10.9.106.11/23", "4", "32", "17")
)

$startRow = 6
$numericCols = @(5, 6, 7)
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        if ($numericCols -contains $col) {
            # Force these to be stored as text (matching the existing
            # sheet's convention of inline-string numeric-looking values)
            # instead of Excel's default numeric auto-detection, then
            # restore the default "Normal" style so no stray formatting
            # is left behind on the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $values[$col - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $values[$col - 1]
        }
    }
}
